$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4, shifting the existing weekly rows (old 4-6) down to (5-7)
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with this week's data, mirroring the static columns
# that repeat for every row of this series (A,B,C,E,F,G,H,I,N,O,Q,R),
# and the new values for the weekly-varying columns (D,J,K,L,M,P).
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 44790
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = 100112012
$ws.Range("G4").Value = "Espinaca"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 8500
$ws.Range("L4").Value = 9000
$ws.Range("M4").Value = 8750
$ws.Range("N4").Value = "`$/cuna 10 kilos"
$ws.Range("O4").Value = "Región Metropolitana"
$ws.Range("P4").Value = 875
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"
